# Added a main workflow that invokes all the other methods
# -> Update the "Test Results" sheet: mark the last Delete test as passed
#    and append the results of the new workflow run (rows 25-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")
$ws.Activate()

# E24 (Delete Test Passed for D09) flips from FALSE to TRUE
$ws.Range("E24").Value = $true

# Seven new fully-tested rows: Read Test Passed (C) and Delete Test Passed (E)
for ($r = 25; $r -le 31; $r++) {
    $ws.Cells.Item($r, 3).Value = $true
    $ws.Cells.Item($r, 5).Value = $true
}

# Two trailing rows only record the Read Test Passed (C) result
$ws.Cells.Item(32, 3).Value = $true
$ws.Cells.Item(33, 3).Value = $true

# Leave the selection where the run ended up
$ws.Range("H28").Select()
